$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns (percentage, multiplier) right after "prices" (column L),
# pushing image name / library / library_base_price / attribute_ids / attribute_texts /
# categories / SEO two columns to the right.
$ws.Range("M1:N1").EntireColumn.Insert()

# New headers
$ws.Range("M1").Value = "percentage"
$ws.Range("N1").Value = "multiplier"

# Row 2 (first data row) gets the sample percentage/multiplier arrays, and keeps its
# original library / library_base_price values of 1 / 9.
$ws.Range("L2").Value = "45;265;550"
$ws.Range("M2").Value = "1,0.8,0.5"
$ws.Range("N2").Value = "1,2,3"
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 9

# Remaining data rows (3-10): percentage/multiplier stay blank, prices string switches to
# semicolons, and library / library_base_price both become 0.
for ($r = 3; $r -le 10; $r++) {
    $ws.Range("L$r").Value = "45;265;550"
    $ws.Range("P$r").Value = 0
    $ws.Range("Q$r").Value = 0
}

# Match the column widths used next to the new columns (best-fit-ish values from the
# original workbook after the insert).
$ws.Columns("M").ColumnWidth = 9.333333333333334
$ws.Columns("N").ColumnWidth = 7.833333333333333
$ws.Columns("R").ColumnWidth = 17.0

# Restore the view state recorded in the saved workbook.
$ws.Range("M3").Select()
